$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "sdsgd"
$ws.Range("B2").Value = "fgfbg"
$ws.Range("C2").Value = "g"
$ws.Range("C3").Value = "b"
$ws.Range("D3").Value = "gfbg"
$ws.Range("G10").Value = "dcdscds"
$ws.Range("H10").Value = "sdcdv"

$ws.Range("H10").Select()
